$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows at row 179 (pushes former rows 179-194 down to 181-196)
$ws.Rows.Item(179).Insert()
$ws.Rows.Item(179).Insert()

# --- New row 179 ---
$ws.Cells.Item(179,1).Value = 7
$ws.Cells.Item(179,2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(179,3).Value = "Ñuble"
$ws.Cells.Item(179,4).Value = 45021
$ws.Cells.Item(179,5).Value = 16
$ws.Cells.Item(179,6).Value = "Fruta"
$ws.Cells.Item(179,7).Value = 100109
$ws.Cells.Item(179,8).Value = "Uva"
$ws.Cells.Item(179,9).Value = 100109001
$ws.Cells.Item(179,10).Value = "Uva"
$ws.Cells.Item(179,11).Value = "Thompson seedless"
$ws.Cells.Item(179,12).Value = "Especial"
$ws.Cells.Item(179,13).Value = 60
$ws.Cells.Item(179,14).Value = 13000
$ws.Cells.Item(179,15).Value = 13000
$ws.Cells.Item(179,16).Value = 13000
$ws.Cells.Item(179,17).Value = "$/bandeja 18 kilos"
$ws.Cells.Item(179,18).Value = "Región de O'Higgins"
$ws.Cells.Item(179,19).Value = 722
$ws.Cells.Item(179,20).Value = 18

# --- New row 180 ---
$ws.Cells.Item(180,1).Value = 7
$ws.Cells.Item(180,2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(180,3).Value = "Ñuble"
$ws.Cells.Item(180,4).Value = 45021
$ws.Cells.Item(180,5).Value = 16
$ws.Cells.Item(180,6).Value = "Fruta"
$ws.Cells.Item(180,7).Value = 100109
$ws.Cells.Item(180,8).Value = "Uva"
$ws.Cells.Item(180,9).Value = 100109001
$ws.Cells.Item(180,10).Value = "Uva"
$ws.Cells.Item(180,11).Value = "Thompson seedless"
$ws.Cells.Item(180,12).Value = "Primera"
$ws.Cells.Item(180,13).Value = 60
$ws.Cells.Item(180,14).Value = 11000
$ws.Cells.Item(180,15).Value = 11000
$ws.Cells.Item(180,16).Value = 11000
$ws.Cells.Item(180,17).Value = "$/bandeja 18 kilos"
$ws.Cells.Item(180,18).Value = "Región de O'Higgins"
$ws.Cells.Item(180,19).Value = 611
$ws.Cells.Item(180,20).Value = 18

# Re-apply date formatting to the new D cells (date-time numeric format), matching the
# rest of the column.
$ws.Range("D179:D180").NumberFormat = "YYYY-MM-DD HH:MM:SS"
